$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.585.70'
$ws.Range("E2").Value = '  -2.63%  '
$ws.Range("D3").Value = '1.656.06'
$ws.Range("E3").Value = '  -4.40%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.10'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -2.44%  '
$ws.Range("E6").Value = '  -2.39%  '
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '24.01'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.59%  '
$ws.Range("E9").Value = '  -2.55%  '
$ws.Range("E10").Value = '  -2.85%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0880'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -1.85%  '
$ws.Range("D12").Value = '1.890.62'
$ws.Range("E12").Value = '  -4.31%  '
$ws.Range("D13").Value = '1.646.92'
$ws.Range("E13").Value = '  -4.80%  '
$ws.Range("E14").Value = '  -2.86%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.564'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +0.00%  '
$ws.Range("E16").Value = '  -2.95%  '
$ws.Range("D17").Value = '27.552.86'
$ws.Range("E17").Value = '  -2.67%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '241.24'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -2.58%  '
$ws.Range("D19").Value = '0.0₃0729'
$ws.Range("E19").Value = '  -3.50%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.55'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -4.75%  '
$ws.Range("E21").Value = '  +0.03%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.47'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -4.05%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.32'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -3.92%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.05'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -2.40%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.79'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -2.39%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.19'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -4.11%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.24'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -3.00%  '
$ws.Range("E28").Value = '  -0.05%  '
$ws.Range("E29").Value = '  -2.49%  '
$ws.Range("E30").Value = '  -0.13%  '
$ws.Range("E31").Value = '  -2.78%  '
$ws.Range("E32").Value = '  -2.97%  '
$ws.Range("D33").Value = '1.452.32'
$ws.Range("E33").Value = '  -2.52%  '
$ws.Range("E34").Value = '  -4.99%  '
$ws.Range("E35").Value = '  -4.98%  '
$ws.Range("E36").Value = '  -0.94%  '
$ws.Range("E37").Value = '  -6.06%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.571'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -5.26%  '
$ws.Range("E39").Value = '  -3.08%  '
$ws.Range("B40").Value = 'WEMIXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.03'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -3.56%  '
$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '69.49'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -1.05%  '
$ws.Range("E42").Value = '  +0.07%  '
$ws.Range("E43").Value = '  -4.07%  '
$ws.Range("E44").Value = '  -3.24%  '
$ws.Range("E45").Value = '  -0.45%  '
$ws.Range("D46").Value = '1.799.22'
$ws.Range("E46").Value = '  -4.23%  '
$ws.Range("E47").Value = '  -1.10%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '88.51'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -2.24%  '
$ws.Range("D49").Value = '0.0₆0106'
$ws.Range("E49").Value = '  -6.22%  '
$ws.Range("E50").Value = '  -1.96%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.84'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -4.27%  '
